$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.3804768025875092
$ws.Range("B1").Value = 3.055658578872681
$ws.Range("C1").Value = 4.62161922454834
$ws.Range("D1").Value = 1.797758817672729
$ws.Range("E1").Value = 0.8158201575279236
